# Auto-generated: apply numeric updates to leve-profit tables across 8 sheets
# matching the authoritative diff (commit: chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 700.53656
$ws.Range("J17").Value = 819.59375
$ws.Range("L17").Value = 2458.78125
$ws.Range("N17").Value = -2794.78125

$ws.Range("H42").Value = 121.8
$ws.Range("I42").Value = 125
$ws.Range("J42").Value = 109
$ws.Range("K42").Value = 375
$ws.Range("L42").Value = 327
$ws.Range("M42").Value = -145
$ws.Range("N42").Value = -787

$ws.Range("H53").Value = 2462.8518
$ws.Range("I53").Value = 2506.6924
$ws.Range("J53").Value = 2422.1428
$ws.Range("K53").Value = 2506.6924
$ws.Range("L53").Value = 2422.1428
$ws.Range("M53").Value = -1869.6924
$ws.Range("N53").Value = -3696.1428

$ws.Range("H132").Value = 1418.0193
$ws.Range("I132").Value = 1173.7709
$ws.Range("K132").Value = 3521.3127
$ws.Range("M132").Value = -991.3126999999999

$ws.Range("H137").Value = 2680.2942
$ws.Range("I137").Value = 2866.4546
$ws.Range("K137").Value = 8599.363799999999
$ws.Range("M137").Value = -6049.363799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3271.5293
$ws.Range("I45").Value = 1848
$ws.Range("K45").Value = 1848
$ws.Range("M45").Value = -1471

$ws.Range("H74").Value = 17132.691
$ws.Range("I74").Value = 21113.375
$ws.Range("J74").Value = 4394.5
$ws.Range("K74").Value = 21113.375
$ws.Range("L74").Value = 4394.5
$ws.Range("M74").Value = -20239.375
$ws.Range("N74").Value = -6142.5

$ws.Range("H77").Value = 17132.691
$ws.Range("I77").Value = 21113.375
$ws.Range("J77").Value = 4394.5
$ws.Range("K77").Value = 105566.875
$ws.Range("L77").Value = 21972.5
$ws.Range("M77").Value = -101198.875
$ws.Range("N77").Value = -30708.5

$ws.Range("H122").Value = 14783.353
$ws.Range("J122").Value = 3861.6667
$ws.Range("L122").Value = 11585.0001
$ws.Range("N122").Value = -16485.0001

$ws.Range("H132").Value = 7289.4546
$ws.Range("I132").Value = 5535.5
$ws.Range("J132").Value = 11966.667
$ws.Range("K132").Value = 16606.5
$ws.Range("L132").Value = 35900.001
$ws.Range("M132").Value = -14076.5
$ws.Range("N132").Value = -40960.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6612.35
$ws.Range("I134").Value = 3461.8333
$ws.Range("K134").Value = 10385.4999
$ws.Range("M134").Value = -7850.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 471.54544
$ws.Range("I22").Value = 454.44446
$ws.Range("K22").Value = 454.44446
$ws.Range("M22").Value = -104.44446

$ws.Range("H31").Value = 9552.683000000001
$ws.Range("I31").Value = 3997.6428
$ws.Range("K31").Value = 3997.6428
$ws.Range("M31").Value = -3702.6428

$ws.Range("H34").Value = 9552.683000000001
$ws.Range("I34").Value = 3997.6428
$ws.Range("K34").Value = 3997.6428
$ws.Range("M34").Value = -3795.6428

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H95").Value = 45975.75
$ws.Range("J95").Value = 45975.75
$ws.Range("L95").Value = 45975.75
$ws.Range("N95").Value = -51467.75

$ws.Range("H122").Value = 1537.5555
$ws.Range("I122").Value = 1511.5
$ws.Range("J122").Value = 1558.4
$ws.Range("K122").Value = 4534.5
$ws.Range("L122").Value = 4675.200000000001
$ws.Range("M122").Value = -2084.5
$ws.Range("N122").Value = -9575.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4317.905
$ws.Range("I34").Value = 182.66667
$ws.Range("J34").Value = 5972
$ws.Range("K34").Value = 548.00001
$ws.Range("L34").Value = 17916
$ws.Range("M34").Value = -464.00001
$ws.Range("N34").Value = -18084

$ws.Range("H39").Value = 6592.143
$ws.Range("I39").Value = 2933.3333
$ws.Range("J39").Value = 7590
$ws.Range("K39").Value = 8799.999899999999
$ws.Range("L39").Value = 22770
$ws.Range("M39").Value = -8505.999899999999
$ws.Range("N39").Value = -23358

$ws.Range("H80").Value = 40004240
$ws.Range("I80").Value = 25003950
$ws.Range("J80").Value = 100005400
$ws.Range("K80").Value = 75011850
$ws.Range("L80").Value = 300016200
$ws.Range("M80").Value = -75010914
$ws.Range("N80").Value = -300018072

$ws.Range("H83").Value = 40004240
$ws.Range("I83").Value = 25003950
$ws.Range("J83").Value = 100005400
$ws.Range("K83").Value = 225035550
$ws.Range("L83").Value = 900048600
$ws.Range("M83").Value = -225030870
$ws.Range("N83").Value = -900057960

$ws.Range("H92").Value = 9617134
$ws.Range("J92").Value = 9617134
$ws.Range("L92").Value = 28851402
$ws.Range("N92").Value = -28853898

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H40").Value = 28000
$ws.Range("J40").Value = 28000
$ws.Range("L40").Value = 28000
$ws.Range("N40").Value = -28302

$ws.Range("H80").Value = 2783.1333
$ws.Range("I80").Value = 2484.3333
$ws.Range("J80").Value = 3231.3333
$ws.Range("K80").Value = 2484.3333
$ws.Range("L80").Value = 3231.3333
$ws.Range("M80").Value = -1486.3333
$ws.Range("N80").Value = -5227.3333

$ws.Range("H83").Value = 2783.1333
$ws.Range("I83").Value = 2484.3333
$ws.Range("J83").Value = 3231.3333
$ws.Range("K83").Value = 12421.6665
$ws.Range("L83").Value = 16156.6665
$ws.Range("M83").Value = -7429.666499999999
$ws.Range("N83").Value = -26140.6665

$ws.Range("H102").Value = 3284.1
$ws.Range("I102").Value = 3146.8696
$ws.Range("J102").Value = 3735
$ws.Range("K102").Value = 3146.8696
$ws.Range("L102").Value = 3735
$ws.Range("M102").Value = -1524.8696
$ws.Range("N102").Value = -6979

$ws.Range("H113").Value = 5764.0244
$ws.Range("I113").Value = 2423.8125
$ws.Range("K113").Value = 2423.8125
$ws.Range("M113").Value = -253.8125

$ws.Range("H122").Value = 3753.2778
$ws.Range("I122").Value = 2007.4445
$ws.Range("J122").Value = 5499.1113
$ws.Range("K122").Value = 6022.333500000001
$ws.Range("L122").Value = 16497.3339
$ws.Range("M122").Value = -3572.333500000001
$ws.Range("N122").Value = -21397.3339

$ws.Range("H126").Value = 3203
$ws.Range("I126").Value = 3203
$ws.Range("K126").Value = 9609
$ws.Range("M126").Value = -7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2073.8
$ws.Range("J22").Value = 2172.6428
$ws.Range("L22").Value = 2172.6428
$ws.Range("N22").Value = -2762.6428

$ws.Range("H23").Value = 23635.334
$ws.Range("I23").Value = 23635.334
$ws.Range("K23").Value = 23635.334
$ws.Range("M23").Value = -23405.334

$ws.Range("H27").Value = 2073.8
$ws.Range("J27").Value = 2172.6428
$ws.Range("L27").Value = 2172.6428
$ws.Range("N27").Value = -2386.6428

$ws.Range("H68").Value = 7928.5713
$ws.Range("I68").Value = 7333.3335
$ws.Range("J68").Value = 8375
$ws.Range("K68").Value = 7333.3335
$ws.Range("L68").Value = 8375
$ws.Range("M68").Value = -6584.3335
$ws.Range("N68").Value = -9873

$ws.Range("H69").Value = 45000
$ws.Range("J69").Value = 45000
$ws.Range("L69").Value = 45000
$ws.Range("N69").Value = -46622

$ws.Range("H71").Value = 7928.5713
$ws.Range("I71").Value = 7333.3335
$ws.Range("J71").Value = 8375
$ws.Range("K71").Value = 36666.6675
$ws.Range("L71").Value = 41875
$ws.Range("M71").Value = -32922.6675
$ws.Range("N71").Value = -49363

$ws.Range("H72").Value = 45000
$ws.Range("J72").Value = 45000
$ws.Range("L72").Value = 135000
$ws.Range("N72").Value = -143112

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H132").Value = 18530142
$ws.Range("I132").Value = 100015800
$ws.Range("J132").Value = 10672.637
$ws.Range("K132").Value = 300047400
$ws.Range("L132").Value = 32017.911
$ws.Range("M132").Value = -300044870
$ws.Range("N132").Value = -37077.911

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 145683.14
$ws.Range("I122").Value = 201656.4
$ws.Range("K122").Value = 604969.2
$ws.Range("M122").Value = -602519.2

$ws.Range("H126").Value = 5856.143
$ws.Range("I126").Value = 2998.3333
$ws.Range("J126").Value = 7999.5
$ws.Range("K126").Value = 8994.999899999999
$ws.Range("L126").Value = 23998.5
$ws.Range("M126").Value = -6524.999899999999
$ws.Range("N126").Value = -28938.5

$ws.Range("H132").Value = 11368539
$ws.Range("I132").Value = 16133802
$ws.Range("K132").Value = 48401406
$ws.Range("M132").Value = -48398876

$ws.Range("H136").Value = 29445564
$ws.Range("I136").Value = 100000800
$ws.Range("K136").Value = 300002400
$ws.Range("M136").Value = -299999850
